# Bugfix: Support px-files without LANGUAGES keyword
# Close #143
#
# The "Table" worksheet holds a keyword/value metadata table (Table3)
# that currently contains a "LANGUAGES" row. Remove that entire row so
# px-files without a LANGUAGES keyword are represented correctly; Excel
# will automatically shift the rows below it up and shrink the table /
# autofilter range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table")

# Locate the row whose first column holds "LANGUAGES" and delete it.
$found = 0
$row = 1
while ($found -eq 0) {
    $key = $ws.Cells.Item($row, 1).Value()
    if ($key -eq "LANGUAGES") {
        $found = $row
    } elseif ($key -eq $null -or $key -eq "") {
        $found = -1
    } else {
        $row = $row + 1
    }
}

if ($found -gt 0) {
    $ws.Rows.Item($found).Delete()
}
